$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45172 = 2023-09-03) for
# every data row (rows 2-471). The commit updates this "last changed" date
# stamp to 45175 (2023-09-06) for all of them.
$ws.Range("C2:C471").Value = 45175
